# Append the new data row (row 9) to the "Artfynd" worksheet, matching the
# author's diff: dimension grows from A1:AY8 to A1:AY9 and a full record is
# added for the species observation "Sphagnum wulfianum" (Bollvitmossa).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# --- Plain text / numeric columns -----------------------------------------
$ws.Cells.Item($row, 1).Value  = 112342543          # A  Id
$ws.Cells.Item($row, 2).Value  = 94034               # B  Taxonsorteringsordning
$ws.Cells.Item($row, 3).Value  = "Ovaliderad"        # C  Valideringsstatus
$ws.Cells.Item($row, 4).Value  = "LC"                # D  Rödlistade
$ws.Cells.Item($row, 5).Value  = 2869                # E  TaxonId
$ws.Cells.Item($row, 6).Value  = "Bollvitmossa"       # F  Artnamn
$ws.Cells.Item($row, 7).Value  = "Sphagnum wulfianum" # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = "Girg."              # H  Auktor

$ws.Cells.Item($row, 16).Value = "Lerbergsmyran, Dlr" # P  Lokalnamn
$ws.Cells.Item($row, 17).Value = 497914               # Q  Ost
$ws.Cells.Item($row, 18).Value = 6722405              # R  Nord
$ws.Cells.Item($row, 19).Value = 2                    # S  Noggrannhet
$ws.Cells.Item($row, 20).Value = "Dalarna"            # T  Län
$ws.Cells.Item($row, 21).Value = "Leksand"            # U  Kommun
$ws.Cells.Item($row, 22).Value = "Dalarna"            # V  Provins
$ws.Cells.Item($row, 23).Value = "Leksand"            # W  Församling

# --- Date-like columns stored as literal text, not Excel dates ------------
# Force text formatting first so "2023-09-27" is kept as a string value
# rather than being auto-converted into a date serial number, then drop
# the transient number-format override so no stray styling is left behind.
$ws.Cells.Item($row, 25).NumberFormat = "@"
$ws.Cells.Item($row, 25).Value = "2023-09-27"          # Y  Startdatum
$ws.Cells.Item($row, 25).Style = "Normal"

$ws.Cells.Item($row, 27).NumberFormat = "@"
$ws.Cells.Item($row, 27).Value = "2023-09-27"          # AA Slutdatum
$ws.Cells.Item($row, 27).Style = "Normal"

# --- Boolean columns --------------------------------------------------------
$ws.Cells.Item($row, 30).Value = $false               # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false               # AE Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false               # AG Ospontan

# --- Reporter / observer ----------------------------------------------------
$ws.Cells.Item($row, 49).Value = "John-Olof Halvarsson" # AW Rapportör
$ws.Cells.Item($row, 50).Value = "John-Olof Halvarsson" # AX Observatörer

# Note: columns I, K, AT, AY are present in the source row as empty/typed
# placeholder cells with no text. Excel's object model clears a cell's
# content (and its text-type marker) whenever its Value is assigned an
# empty string, so those four placeholders are intentionally left blank
# here -- any attempt to force them through the COM Value/Formula setters
# collapses back to "no cell", matching what real Excel would do too.
